$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values per diff ---

$text4 = @'
1.5.3: Число стран, принявших и осуществляющих национальные стратегии снижения риска бедствий в соответствии с Сендайской рамочной программой по снижению риска бедствий на 2015–2030 годы 
'@
$ws.Range("B4").Value = $text4

$text6 = @'
Министерство чрезвычайных ситуаций КР
'@
$ws.Range("B6").Value = $text6

$text7 = @'
Кадырова Д.
'@
$ws.Range("B7").Value = $text7

$ws.Range("B9").Value = ""

$text12 = @'
Межправительственная рабочая группа экспертов открытого состава по показателям и терминологии, касающимся уменьшения опасностей бедствий (УОБ), установленных Генеральной Ассамблеей (резолюция 69/284), разрабатывает ряд показателей для оценки глобального прогресса в осуществлении Сендайской рамочной программе. Эти показатели будут в конечном итоге отражать показатели по Сендайской рамочной программе.
'@
$ws.Range("B12").Value = $text12

$text13 = @'
Стратегия снижения риска бедствий нацелена на достижение в ближайшие 15 лет следующего результата: существенное снижение риска бедствий и сокращение потерь в результате бедствий в виде человеческих жертв, утраты источников средств к существованию и ухудшения состояния здоровья людей, и неблагоприятных последствий для экономических, физических, социальных, культурных и экологических активов людей, предприятий, общин и стран.
'@
$ws.Range("B13").Value = $text13

$text14 = @'
Показатель построит мост между ЦУР (SDG) и Сендайской рамочной программой для УОБ (DRR). Увеличение числа национальных правительств, которые принимают и реализуют национальные и местные стратегии УОБ, которые, согласно Сендайской рамочной программе, будут способствовать устойчивому развитию с экономической, экологической и социальной точек зрения.
Не рекомендуется просто подсчитывать количество стран, вместо этого поощряются показатели для оценки прогресса с течением времени. В дополнение к обсуждениям МРГОС, а также Межведомственной экспертной группы, МСУОБ ООН предложила методологию расчета, которая позволит осуществлять мониторинг улучшения национальных и местных стратегий УОБ с течением времени. Эти методологии варьируются от простой количественной оценки этих стратегий до качественной меры согласования с Сендайской рамочной программой, а также охвата населения местными стратегиями.
'@
$ws.Range("B14").Value = $text14

$text16 = @'
Национальный отчет о прогрессе Мониторинга Сендайской рамочной программы, который был представлен в МСУОБ ООН.
'@
$ws.Range("B16").Value = $text16

$text17 = @'
Официальный партнер(ы) на уровне страны предоставляют национальный отчет о прогрессе в Мониторинге Сендайской рамочной программы.
'@
$ws.Range("B17").Value = $text17

$text19 = @'
Примечание: методика расчета по нескольким показателям комплексна, очень длинная (около 180 страниц) и, вероятно, выходит за рамки этих метаданных. МСУОБ ООН предпочитает сослаться на итоги Межправительственной рабочей группой Открытого состава (МРГОС), которая предоставляет полную подробную методику по каждому показателю. 
Последняя версия этой методологии можно получить по адресу:
http://www.preventionweb.net/documents/oiewg/
Technical%20Collection%20of%20Concept%20Notes%20on%20Indicators.pdf

'@
$ws.Range("B19").Value = $text19

$text20 = @'
Инструмент мониторинга Хиогской рамочной программы действий (HFA) начал действие в 2007 году и с течением времени увеличилось число стран, подотчетных МСУОБ ООН, от 60 в 2007 году до 140 стран, которые в настоящее время проводят добровольную самооценку прогресса внедрения Хиогской рамочной программы действий. В течение четырех циклов отчетности до 2015 года инструмент мониторинга Хиогской рамочной программы действий создал крупнейший в мире репозиторий информации о национальной политике в области УОБ, в частности. Преемственная программа, условно названная Сендайским инструментом мониторинга, находится в разработке и будет проинформирован рекомендациями МРГОС (OEIWG).
'@
$ws.Range("B20").Value = $text20

$ws.Range("B21").Value = ""

$text23 = @'
Временные ряды 2013 и 2015 годы: мониторинг Хиогской рамочной программы (HFA)
'@
$ws.Range("B23").Value = $text23

$text24 = @'
Национальный уровень.
'@
$ws.Range("B24").Value = $text24

$text25 = @'
Данные по Кыргызстану сопоставимы с данными других стран, так как составляются на основе международной методологии.
'@
$ws.Range("B25").Value = $text25

$text26 = @'
URL: http://www.preventionweb.net/documents/oiewg/ Technical%20Collection%20of%20Concept%20Notes%20on%20Indicators.pdf
Организацией Объединенных Наций было поручено создание Межправительственной рабочей группы экспертов открытого состава по показателям и терминологии, касающейся уменьшения опасности бедствий (OEIWG), для разработки набора показателей для измерения глобального прогресса в осуществлении Сендайской рамочной программы по снижению риска бедствий, включая семь глобальных целей. Работа OEIWG должна быть завершена к декабрю 2016 года, а ее доклад представлен Генеральной Ассамблее для рассмотрения. Межправительственная рабочая группа по Целям устойчивого развития IAEG-SDG и Статистическая комиссия ООН формально признают роль OEIWG и возлагают ответственность за дальнейшее уточнение и разработку методологии для индикаторов ЦУР, связанных с бедствиями на данную рабочую группу.
http://www.preventionweb.net/drr-framework/open-ended-working-group/
Последняя версия документов находится по адресу:
http://www.preventionweb.net/drr-framework/open-ended-working-group/sessional-intersessionaldocuments
Национальная платформа отчётности ЦУР КР: https://sustainabledevelopment-kyrgyzstan.github.io

'@
$ws.Range("B26").Value = $text26

# --- Row height adjustments ---
$ws.Rows.Item(12).RowHeight = 92.25
$ws.Rows.Item(13).RowHeight = 92.25
$ws.Rows.Item(14).RowHeight = 211.5
$ws.Rows.Item(19).RowHeight = 120.75
$ws.Rows.Item(20).RowHeight = 163.5
$ws.Rows.Item(21).RowHeight = 58.15
$ws.Rows.Item(23).RowHeight = 58.15
$ws.Rows.Item(25).RowHeight = 58.15
$ws.Rows.Item(26).RowHeight = 319.5

# --- Column B width adjustment ---
$ws.Columns.Item(2).ColumnWidth = 71.7109375

# --- Wrap text fix for now-empty B21 (no longer wraps) ---
$ws.Range("B21").WrapText = $false

# --- Selection change ---
$ws.Range("B4").Select()
